$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "verbs"

$ws.Range('A1').Value = 'blast'
$ws.Range('B1').Value = 597
$ws.Range('A2').Value = 'break'
$ws.Range('B2').Value = 564
$ws.Range('A3').Value = 'brief'
$ws.Range('B3').Value = 574
$ws.Range('A4').Value = 'change'
$ws.Range('B4').Value = 563
$ws.Range('A5').Value = 'climb'
$ws.Range('B5').Value = 598
$ws.Range('A6').Value = 'close'
$ws.Range('B6').Value = 538
$ws.Range('A7').Value = 'cross'
$ws.Range('B7').Value = 532
$ws.Range('A8').Value = 'cut'
$ws.Range('B8').Value = 544
$ws.Range('A9').Value = 'dig'
$ws.Range('B9').Value = 582
$ws.Range('A10').Value = 'drink'
$ws.Range('B10').Value = 555
$ws.Range('A11').Value = 'drop'
$ws.Range('B11').Value = 535
$ws.Range('C11').Value = 'free, dr'
$ws.Range('A12').Value = 'eat'
$ws.Range('B12').Value = 554
$ws.Range('A13').Value = 'enter'
$ws.Range('B13').Value = 526
$ws.Range('C13').Value = 'in'
$ws.Range('A14').Value = 'exit'
$ws.Range('B14').Value = 527
$ws.Range('C14').Value = 'out'
$ws.Range('A15').Value = 'extinguish'
$ws.Range('B15').Value = 546
$ws.Range('C15').Value = 'off'
$ws.Range('A16').Value = 'feed'
$ws.Range('B16').Value = 559
$ws.Range('A17').Value = 'fill'
$ws.Range('B17').Value = 560
$ws.Range('A18').Value = 'go'
$ws.Range('B18').Value = 528
$ws.Range('C18').Value = 'walk, run'
$ws.Range('A19').Value = 'get'
$ws.Range('B19').Value = 534
$ws.Range('C19').Value = 'g, take'
$ws.Range('A20').Value = 'help'
$ws.Range('B20').Value = 580
$ws.Range('C20').Value = '?'
$ws.Range('A21').Value = 'history'
$ws.Range('B21').Value = 576
$ws.Range('A22').Value = 'hit'
$ws.Range('B22').Value = 550
$ws.Range('A23').Value = 'hoot'
$ws.Range('B23').Value = 592
$ws.Range('A24').Value = 'info'
$ws.Range('B24').Value = 585
$ws.Range('A25').Value = 'inventory'
$ws.Range('B25').Value = 594
$ws.Range('C25').Value = 'I '
$ws.Range('A26').Value = 'jump'
$ws.Range('B26').Value = 531
$ws.Range('A27').Value = 'kill'
$ws.Range('B27').Value = 549
$ws.Range('A28').Value = 'lift'
$ws.Range('B28').Value = 561
$ws.Range('A29').Value = 'light'
$ws.Range('B29').Value = 545
$ws.Range('C29').Value = 'on'
$ws.Range('A30').Value = 'listen'
$ws.Range('B30').Value = 541
$ws.Range('A31').Value = 'look'
$ws.Range('B31').Value = 571
$ws.Range('C31').Value = 'l, examine, x'
$ws.Range('A32').Value = 'news'
$ws.Range('B32').Value = 557
$ws.Range('A33').Value = 'open'
$ws.Range('B33').Value = 536
$ws.Range('C33').Value = 'unlock'
$ws.Range('A34').Value = 'play'
$ws.Range('B34').Value = 552
$ws.Range('A35').Value = 'pour'
$ws.Range('B35').Value = 553
$ws.Range('A36').Value = 'quit'
$ws.Range('B36').Value = 572
$ws.Range('C36').Value = 'q'
$ws.Range('A37').Value = 'read'
$ws.Range('B37').Value = 562
$ws.Range('A38').Value = 'reflect'
$ws.Range('B38').Value = 603
$ws.Range('A39').Value = 'remove'
$ws.Range('B39').Value = 540
$ws.Range('A40').Value = 'rest'
$ws.Range('B40').Value = 491
$ws.Range('A41').Value = 'restore'
$ws.Range('B41').Value = 490
$ws.Range('A42').Value = 'retreat'
$ws.Range('B42').Value = 530
$ws.Range('C42').Value = 'back'
$ws.Range('A43').Value = 'ride'
$ws.Range('B43').Value = 565
$ws.Range('A44').Value = 'rub'
$ws.Range('B44').Value = 556
$ws.Range('C44').Value = 'polish'
$ws.Range('A45').Value = 'save'
$ws.Range('B45').Value = 489
$ws.Range('A46').Value = 'say'
$ws.Range('B46').Value = 493
$ws.Range('A47').Value = 'score'
$ws.Range('B47').Value = 595
$ws.Range('C47').Value = 'sc'
$ws.Range('A48').Value = 'scry'
$ws.Range('B48').Value = 604
$ws.Range('A49').Value = 'tame'
$ws.Range('B49').Value = 548
$ws.Range('C49').Value = 'placate'
$ws.Range('A50').Value = 'terse'
$ws.Range('B50').Value = 573
$ws.Range('A51').Value = 'throw'
$ws.Range('B51').Value = 557
$ws.Range('C51').Value = 'th'
$ws.Range('A52').Value = 'translate'
$ws.Range('B52').Value = 542
$ws.Range('A53').Value = 'verbose'
$ws.Range('B53').Value = 575
$ws.Range('A54').Value = 'wash'
$ws.Range('B54').Value = 610
$ws.Range('A55').Value = 'wave'
$ws.Range('B55').Value = 547
$ws.Range('C55').Value = 'swing'
$ws.Range('A56').Value = 'wear'
$ws.Range('B56').Value = 539
